# "chup man hinh update test loi" - add screenshot/test-error rows 18-24
# (sheet rows 20-26) to the bug tracker sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 20 (STT 18) ---
$ws.Range("A20").Value = 18
$ws.Range("B18").Copy()
[void]$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").HorizontalAlignment = -4108
$ws.Range("B20").Value = "không cho tùy chỉnh loại khách hàng lúc thêm, edit"
[void]$ws.Range("B18").Copy()
[void]$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("E20").Value = "lúc thêm mặc định là khách lẻ"

# --- Row 21 (STT 19) ---
$ws.Range("A21").Value = 19
[void]$ws.Range("A20").Copy()
[void]$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("B21").Value = "header chi tiết phiếu nhập nên để là mã chứng từ thay cho id"
[void]$ws.Range("B20").Copy()
[void]$ws.Range("B21").PasteSpecial(-4122)

# --- Row 22 (STT 20) ---
$ws.Range("A22").Value = 20
[void]$ws.Range("A21").Copy()
[void]$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("B22").Value = "trong chi tiet phieu nhap thieu thong tin nha cung cap"
[void]$ws.Range("B21").Copy()
[void]$ws.Range("B22").PasteSpecial(-4122)

# --- Row 23 (STT 21) ---
$ws.Range("A23").Value = 21
[void]$ws.Range("A22").Copy()
[void]$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("B23").Value = "xuất file excel danh sách chi nhánh trong module nhập xuất chưa hoàn chỉnh"
[void]$ws.Range("B22").Copy()
[void]$ws.Range("B23").PasteSpecial(-4122)

# --- Row 24 (STT 22) ---
$ws.Range("A24").Value = 22
[void]$ws.Range("A23").Copy()
[void]$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("B24").Value = "báo cáo doanh thu theo chi nhánh: phần chọn tất cả chi nhánh hình như chưa làm"
[void]$ws.Range("B23").Copy()
[void]$ws.Range("B24").PasteSpecial(-4122)

# --- Row 25 (STT 23) ---
$ws.Range("A25").Value = 23
[void]$ws.Range("A24").Copy()
[void]$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("E25").Value = "nếu rảnh "
$ws.Range("B25").Value = "báo cáo: chỉnh header thể hiện biễu đồ đang biễu diễn về sản phẩm j, chi nhánh nào…"
[void]$ws.Range("B24").Copy()
[void]$ws.Range("B25").PasteSpecial(-4122)

# --- Row 26 (STT 24) ---
$ws.Range("A26").Value = 24
[void]$ws.Range("A25").Copy()
[void]$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("B26").Value = "tắt hiển thị phần thao tác trong phân quyền"
[void]$ws.Range("B25").Copy()
[void]$ws.Range("B26").PasteSpecial(-4122)

# Widen column B to fit the new (longer) text, matching the bestFit recalculation.
$ws.Columns.Item(2).ColumnWidth = 73.8

# Move the active cell below the new last row, like the original author did
# after typing the final entry.
[void]$ws.Range("B27").Select()
